$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Seed Surveys")

# Fill in the R column ("Agrin_Endo") values for rows 2-51 that were
# previously blank (style-only) cells, matching the author's recorded data.
$values = @{
    2  = 1
    3  = 1
    4  = 1
    5  = 1
    6  = 0
    7  = 1
    8  = 0
    9  = 1
    10 = 1
    11 = 1
    12 = 0
    13 = 0
    16 = 1
    17 = 1
    18 = 1
    19 = 1
    20 = 0
    21 = 0
    22 = 0
    23 = 1
    24 = 0
    25 = 1
    26 = 0
    27 = 0
    30 = 0
    31 = 1
    32 = 1
    33 = 0
    34 = 1
    35 = 1
    36 = 0
    37 = 0
    38 = 1
    39 = 1
    40 = 0
    41 = 1
    44 = 0
    45 = 1
    46 = 1
    47 = 0
    48 = 1
    49 = 0
    50 = 0
    51 = 1
}

foreach ($row in $values.Keys) {
    $ws.Range("R$row").Value = $values[$row]
}

# Row 43's R cell stays conceptually blank, but gets fully cleared
# (content + formatting) rather than left as a styled-empty cell.
$ws.Range("R43").Clear()

# Restore view state: scroll position and active selection on the sheet.
$excel.ActiveWindow.ScrollRow = 19
$ws.Range("R52").Select()
